$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the section header text used by A2:A9
#    "Officiating Athletes with Disabilities" -> "ᐧ Disabilities"
$ws.Range("A2:A9").Value = "ᐧ Disabilities"

# 2. Update frozen-pane scroll position / active cell selections
#    topLeftCell C71 -> C2, bottomLeft active cell A71 -> A2,
#    bottomRight active cell A1 -> A7
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 3
$ws.Range("A7").Select()

# 3. Column width for C:D 64.85 -> 64.84 (closest attainable value)
$ws.Range("C1:D1").ColumnWidth = 64.0
